$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.079.42'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").Value = '2.413.62'
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '554.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.79'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +0.72%  '
$ws.Range("E9").Value = '  -0.97%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.68'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.47%  '
$ws.Range("E11").Value = '  -0.66%  '
$ws.Range("E12").Value = '  -1.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.80'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.07%  '
$ws.Range("D14").Value = '2.845.40'
$ws.Range("E14").Value = '  -0.15%  '
$ws.Range("D15").Value = '59.981.75'
$ws.Range("E15").Value = '  +0.09%  '
$ws.Range("E16").Value = '  -0.23%  '
$ws.Range("D17").Value = '2.411.90'
$ws.Range("E17").Value = '  -0.75%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.22'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.84%  '
$ws.Range("E19").Value = '  +3.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '326.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.79'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.95%  '
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '64.77'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.179'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.63'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("E26").Value = '  +0.32%  '
$ws.Range("E27").Value = '  +4.92%  '
$ws.Range("E28").Value = '  -1.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.77'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.76'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.73%  '
$ws.Range("E32").Value = '  +5.22%  '
$ws.Range("E33").Value = '  -3.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.45'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.94%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("E36").Value = '  +2.14%  '
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("E38").Value = '  +0.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '324.52'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.89%  '
$ws.Range("E40").Value = '  -1.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '146.17'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.97%  '
$ws.Range("E42").Value = '  -1.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0963'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.75'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.63%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0515'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.577'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("E47").Value = '  -1.32%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '11.05'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("E49").Value = '  -0.96%  '
$ws.Range("E51").Value = '  -1.65%  '
